$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters: longer/more-specific tokens first to avoid collisions
# (e.g. D80 must not get touched by a later D8x rule, etc.)
$used = $ws.UsedRange

$used.Replace("D80", "D86", -4163)
$used.Replace("D64", "D69", -4163)
$used.Replace("D51", "D55", -4163)
$used.Replace("S30", "S31", -4163)
